$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("drivers_data")

# New drivers being added to the bottom of the table (matches the row
# directly above it for formatting: Name / Tel / Address / Yes,Yes,Yes / No,No,No)
$newDrivers = @(
    @("GMC BREAKDOWN RECOVERY LTD", "07838 666656", "39 Chequers Way, Woodley, Reading, England, RG5 3EH"),
    @("ROAD 2 RECOVERY", "07479 032059", "27 Hale Ln, Otford, Sevenoaks TN14 5NP"),
    @("247 TYRE MOBILE LTD", "020 3488 6465", "71-75 Shelton Street, Covent Garden WC2H 9JQ")
)

# Last populated row currently in the table (row 53 before this edit)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($i = 0; $i -lt $newDrivers.Count; $i++) {
    $srcRow = $lastRow + $i
    $destRow = $lastRow + $i + 1
    $data = $newDrivers[$i]

    # Clone the row immediately above (formatting + row height) onto the new row,
    # then overwrite the Name/Tel/Address cells and the Yes/No service flags.
    $ws.Range("A" + $srcRow + ":I" + $srcRow).Copy($ws.Range("A" + $destRow + ":I" + $destRow))

    $ws.Cells.Item($destRow, 1).Value = $data[0]
    $ws.Cells.Item($destRow, 2).Value = $data[1]
    $ws.Cells.Item($destRow, 3).Value = $data[2]
    $ws.Cells.Item($destRow, 4).Value = "Yes"
    $ws.Cells.Item($destRow, 5).Value = "Yes"
    $ws.Cells.Item($destRow, 6).Value = "Yes"
    $ws.Cells.Item($destRow, 7).Value = "No"
    $ws.Cells.Item($destRow, 8).Value = "No"
    $ws.Cells.Item($destRow, 9).Value = "No"
}
$excel.CutCopyMode = $false

# Scroll the view down and leave the selection where editing finished
$ws.Activate()
$ws.Range("C" + $lastRow).Select()
